# Update localization status report: files "46d5015d-472f-4702-9b37-e4cdd99e6287.md"
# and "b6b08c00-298a-4d88-91da-9e2e70996fdf.md" have moved from "Ready for handoff"
# to "In Translation" for both the zh-cn and de-de locales. Reflect this on the
# Overview roll-up sheet as well as the per-locale detail sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: rows 3 & 4 (46d5015d..., b6b08c00...) -> zh-cn (B) / de-de (C) columns
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "In Translation"
$wsOverview.Range("C3").Value = "In Translation"
$wsOverview.Range("B4").Value = "In Translation"
$wsOverview.Range("C4").Value = "In Translation"

# zh-cn detail sheet: Status column (C) for rows 3 & 4
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Range("C4").Value = "In Translation"

# de-de detail sheet: Status column (C) for rows 3 & 4
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Range("C4").Value = "In Translation"
